# Apply updated dSF (column F) values as part of a data re-pull / mean
# calculation refresh. Only the rows whose computed dSF value changed
# after the repull are touched; all other data is left untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    2  = -1
    3  = 1
    4  = -1
    5  = -3
    6  = -1
    8  = -1
    9  = 3
    10 = -1
    12 = 1
    13 = 2
    15 = 4
    16 = 4
    17 = 3
    18 = 2
    19 = 7
    21 = -2
    22 = 1
    25 = 11
}

foreach ($row in $updates.Keys) {
    $ws.Range("F$row").Value = $updates[$row]
}
